$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.537.66'
$ws.Range('E2').Value = '  +5.58%  '
$ws.Range('D3').Value = '1.722.28'
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '225.51'
$ws.Range('E5').Value = '  +3.38%  '
$ws.Range('D6').Value = '0.5358'
$ws.Range('E6').Value = '  +3.09%  '
$ws.Range('D7').Value = '1.004'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '0.2668'
$ws.Range('E8').Value = '  +1.08%  '
$ws.Range('D9').Value = '0.06596'
$ws.Range('E9').Value = '  +4.22%  '
$ws.Range('D10').Value = '21.70'
$ws.Range('E10').Value = '  +6.63%  '
$ws.Range('D11').Value = '0.07713'
$ws.Range('E11').Value = '  +0.41%  '
$ws.Range('D12').Value = '4.615'
$ws.Range('E12').Value = '  +0.47%  '
$ws.Range('D13').Value = '1.729.67'
$ws.Range('E13').Value = '  +4.52%  '
$ws.Range('D14').Value = '1.960.32'
$ws.Range('E14').Value = '  +4.35%  '
$ws.Range('D15').Value = '0.5844'
$ws.Range('E15').Value = '  +4.72%  '
$ws.Range('D16').Value = '0.0₅8305'
$ws.Range('E16').Value = '  +2.23%  '
$ws.Range('D17').Value = '67.91'
$ws.Range('E17').Value = '  +3.93%  '
$ws.Range('D18').Value = '27.550.50'
$ws.Range('E18').Value = '  +5.67%  '
$ws.Range('D19').Value = '220.36'
$ws.Range('E19').Value = '  +15.31%  '
$ws.Range('E20').Value = '  +0.01%  '
$ws.Range('D21').Value = '4.722'
$ws.Range('E21').Value = '  +2.13%  '
$ws.Range('E22').Value = '  +1.51%  '
$ws.Range('D23').Value = '6.088'
$ws.Range('E23').Value = '  +2.94%  '
$ws.Range('E24').Value = '  +0.05%  '
$ws.Range('D25').Value = '148.48'
$ws.Range('E25').Value = '  +2.97%  '
$ws.Range('D26').Value = '1.742'
$ws.Range('E26').Value = '  +16.04%  '
$ws.Range('E27').Value = '  +4.23%  '
$ws.Range('D28').Value = '7.407'
$ws.Range('E28').Value = '  +2.77%  '
$ws.Range('D29').Value = '16.65'
$ws.Range('E29').Value = '  +4.70%  '
$ws.Range('D30').Value = '0.05600'
$ws.Range('E30').Value = '  +2.44%  '
$ws.Range('E31').Value = '  +2.74%  '
$ws.Range('E32').Value = '  +3.42%  '
$ws.Range('D33').Value = '3.452'
$ws.Range('E33').Value = '  +3.10%  '
$ws.Range('E34').Value = '  +6.77%  '
$ws.Range('D35').Value = '2.833'
$ws.Range('E35').Value = '  +1.67%  '
$ws.Range('D36').Value = '0.9611'
$ws.Range('E36').Value = '  +1.68%  '
$ws.Range('D37').Value = '2.429'
$ws.Range('E37').Value = '  +0.30%  '
$ws.Range('D38').Value = '0.5959'
$ws.Range('E38').Value = '  +5.59%  '
$ws.Range('D39').Value = '0.01654'
$ws.Range('E39').Value = '  +4.90%  '
$ws.Range('D40').Value = '5.923'
$ws.Range('E40').Value = '  +1.30%  '
$ws.Range('D41').Value = '0.8530'
$ws.Range('E41').Value = '  +3.25%  '
$ws.Range('D42').Value = '1.053.77'
$ws.Range('E42').Value = '  +2.71%  '
$ws.Range('D43').Value = '1.004'
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('D44').Value = '101.27'
$ws.Range('E44').Value = '  +0.20%  '
$ws.Range('D45').Value = '1.866.46'
$ws.Range('E45').Value = '  +4.27%  '
$ws.Range('E46').Value = '  +3.72%  '
$ws.Range('D47').Value = '59.08'
$ws.Range('E47').Value = '  +2.77%  '
$ws.Range('D48').Value = '8.195'
$ws.Range('E48').Value = '  +2.53%  '
$ws.Range('D49').Value = '0.4436'
$ws.Range('E49').Value = '  +2.37%  '
$ws.Range('D50').Value = '1.000'
$ws.Range('E50').Value = '  +0.12%  '
$ws.Range('D51').Value = '0.05252'
$ws.Range('E51').Value = '  +1.65%  '
